$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.111.71'
$ws.Range("E2").Value = '  -2.85%  '
$ws.Range("D3").Value = '1.652.13'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.51'
$ws.Range("E5").Value = '  -1.87%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4783'
$ws.Range("E7").Value = '  -8.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2611'
$ws.Range("E8").Value = '  -4.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.05959'
$ws.Range("E9").Value = '  -3.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07066'
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("D11").Value = '1.660.67'
$ws.Range("E11").Value = '  -4.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.39'
$ws.Range("E12").Value = '  -3.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6177'
$ws.Range("E13").Value = '  -3.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.577'
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '72.85'
$ws.Range("E15").Value = '  -5.89%  '
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").Value = '25.108.95'
$ws.Range("E18").Value = '  -2.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.35'
$ws.Range("E19").Value = '  -3.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006518'
$ws.Range("E20").Value = '  -3.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.406'
$ws.Range("E21").Value = '  +3.19%  '
$ws.Range("D22").Value = '1.869.21'
$ws.Range("E22").Value = '  -4.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.445'
$ws.Range("E23").Value = '  -2.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.258'
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '132.78'
$ws.Range("E25").Value = '  -4.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.71'
$ws.Range("E26").Value = '  -3.14%  '
$ws.Range("E27").Value = '  -8.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.704'
$ws.Range("E28").Value = '  -3.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '101.90'
$ws.Range("E29").Value = '  -3.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.808'
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07885'
$ws.Range("E31").Value = '  -4.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.516'
$ws.Range("E32").Value = '  -4.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04590'
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.606'
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9423'
$ws.Range("E35").Value = '  -4.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.5827'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.621'
$ws.Range("E37").Value = '  -2.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01535'
$ws.Range("E38").Value = '  -4.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.002'
$ws.Range("E39").Value = '  +0.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8359'
$ws.Range("E40").Value = '  +12.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.832'
$ws.Range("E41").Value = '  -5.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.76'
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3687'
$ws.Range("E43").Value = '  -4.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.826'
$ws.Range("E44").Value = '  -3.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1125'
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.027'
$ws.Range("E46").Value = '  -3.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05143'
$ws.Range("E47").Value = '  -1.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '51.94'
$ws.Range("E48").Value = '  -4.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.41'
$ws.Range("E49").Value = '  -3.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.002'
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.304'
$ws.Range("E51").Value = '  -3.54%  '
